$d = $word.ActiveDocument

# The heading paragraph reads:
#   "Ana Caroline, Débora Loro, Matheus e Natália - 423"
# and the hyphen-minus separating the authors from the class number
# must become an en dash ("–"), with the text split into three runs
# (the text before the dash, the dash itself, and the text after it)
# just like Word does when you select only the dash and retype it.
$dash = $d.Content
$found = $dash.Find.Execute("-", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "–", 2)

if (-not $found) {
    throw "Could not find the hyphen-minus to replace in the document."
}

# $dash now covers exactly the newly-inserted en dash. Re-apply the
# paragraph's explicit run formatting (Times New Roman, bold) to it so
# it keeps matching the surrounding runs' look, which is what causes
# Word to keep it as its own, identically-formatted run rather than
# silently merging back into the neighbouring text.
$dash.Font.Name = "Times New Roman"
$dash.Font.NameBi = "Times New Roman"
$dash.Bold = 1
